# Fix "Recorded By" (column G) entries where the literal token "System"
# is listed first in a comma-separated list of recorders. The system
# account should be listed last instead of first, so swap the first and
# last comma-separated tokens whenever the first token is exactly "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1 -and $parts[0] -ceq "System") {
        $first = $parts[0]
        $last = $parts[$parts.Length - 1]
        $parts[0] = $last
        $parts[$parts.Length - 1] = $first

        $newText = [string]::Join(", ", $parts)
        $cell.Value2 = $newText
    }
}
